$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Update the summary sheet ("总计"): insert a new 2022-Q3 row
#    at the top of the data table, shifting existing quarters down
#    by one row. The index column (A) keeps its original sequence
#    of values per row position, it is not shifted.
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Capture the last existing row (old row 7, 2021-Q1) into the new
# row 8 first -- before it gets overwritten by the cascading shift
$summary.Range("A7:D7").Copy($summary.Range("A8:D8"))
$summary.Range("A8").Value = 6

# Shift B:D values down one row at a time, bottom-up, carrying the
# (lack of) formatting along with Copy so no extra styles appear
$summary.Range("B6:D6").Copy($summary.Range("B7:D7"))
$summary.Range("B5:D5").Copy($summary.Range("B6:D6"))
$summary.Range("B4:D4").Copy($summary.Range("B5:D5"))
$summary.Range("B3:D3").Copy($summary.Range("B4:D4"))
$summary.Range("B2:D2").Copy($summary.Range("B3:D3"))

# Write the new 2022-Q3 figures into row 2
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 21
$summary.Range("D2").Value = 2.72

# ---------------------------------------------------------------
# 2) Insert a brand-new "2022-Q3" worksheet right after "总计"
#    (i.e. right before the current "2022-Q2" sheet), cloned from
#    the "2022-Q2" sheet so it keeps identical formatting, then
#    overwrite its values with the Q3 fund holdings.
# ---------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q2Sheet.Copy($q2Sheet, $null)
$q3Sheet = $wb.Worksheets.Item($q2Sheet.Index - 1)
$q3Sheet.Name = "2022-Q3"

# Force a Text number format on the numeric-looking text columns so
# that values such as "002621" or "10.65" are kept as text instead
# of being coerced into numbers (which would also lose leading
# zeroes on fund codes)
$q3Sheet.Range("B2:B22").NumberFormat = "@"
$q3Sheet.Range("D2:G22").NumberFormat = "@"

$q3Sheet.Range("A2").Value = 0
$q3Sheet.Range("B2").Value = '002621'
$q3Sheet.Range("C2").Value = '中欧消费主题股票A'
$q3Sheet.Range("D2").Value = '10.65'
$q3Sheet.Range("E2").Value = '88.44'
$q3Sheet.Range("F2").Value = '6.23'
$q3Sheet.Range("G2").Value = '0.6635'
$q3Sheet.Range("H2").Value = 6
$q3Sheet.Range("A3").Value = 1
$q3Sheet.Range("B3").Value = '003751'
$q3Sheet.Range("C3").Value = '万家瑞隆混合A'
$q3Sheet.Range("D3").Value = '14.54'
$q3Sheet.Range("E3").Value = '93.96'
$q3Sheet.Range("F3").Value = '2.85'
$q3Sheet.Range("G3").Value = '0.4144'
$q3Sheet.Range("H3").Value = 8
$q3Sheet.Range("A4").Value = 2
$q3Sheet.Range("B4").Value = '006972'
$q3Sheet.Range("C4").Value = '金鹰民安回报一年定期开放混合A'
$q3Sheet.Range("D4").Value = '26.02'
$q3Sheet.Range("E4").Value = '22.96'
$q3Sheet.Range("F4").Value = '1.55'
$q3Sheet.Range("G4").Value = '0.4033'
$q3Sheet.Range("H4").Value = 4
$q3Sheet.Range("A5").Value = 3
$q3Sheet.Range("B5").Value = '002697'
$q3Sheet.Range("C5").Value = '中欧消费主题股票C'
$q3Sheet.Range("D5").Value = '5.37'
$q3Sheet.Range("E5").Value = '88.44'
$q3Sheet.Range("F5").Value = '6.23'
$q3Sheet.Range("G5").Value = '0.3346'
$q3Sheet.Range("H5").Value = 6
$q3Sheet.Range("A6").Value = 4
$q3Sheet.Range("B6").Value = '010852'
$q3Sheet.Range("C6").Value = '中欧内需成长混合A'
$q3Sheet.Range("D6").Value = '3.70'
$q3Sheet.Range("E6").Value = '90.11'
$q3Sheet.Range("F6").Value = '6.40'
$q3Sheet.Range("G6").Value = '0.2368'
$q3Sheet.Range("H6").Value = 5
$q3Sheet.Range("A7").Value = 5
$q3Sheet.Range("B7").Value = '005620'
$q3Sheet.Range("C7").Value = '中欧品质消费股票A'
$q3Sheet.Range("D7").Value = '1.62'
$q3Sheet.Range("E7").Value = '91.58'
$q3Sheet.Range("F7").Value = '6.12'
$q3Sheet.Range("G7").Value = '0.0991'
$q3Sheet.Range("H7").Value = 6
$q3Sheet.Range("A8").Value = 6
$q3Sheet.Range("B8").Value = '519125'
$q3Sheet.Range("C8").Value = '浦银安盛消费升级混合A'
$q3Sheet.Range("D8").Value = '2.17'
$q3Sheet.Range("E8").Value = '90.82'
$q3Sheet.Range("F8").Value = '4.19'
$q3Sheet.Range("G8").Value = '0.0909'
$q3Sheet.Range("H8").Value = 10
$q3Sheet.Range("A9").Value = 7
$q3Sheet.Range("B9").Value = '005621'
$q3Sheet.Range("C9").Value = '中欧品质消费股票C'
$q3Sheet.Range("D9").Value = '0.97'
$q3Sheet.Range("E9").Value = '91.58'
$q3Sheet.Range("F9").Value = '6.12'
$q3Sheet.Range("G9").Value = '0.0594'
$q3Sheet.Range("H9").Value = 6
$q3Sheet.Range("A10").Value = 8
$q3Sheet.Range("B10").Value = '004818'
$q3Sheet.Range("C10").Value = '国寿安保目标策略灵活配置混合A'
$q3Sheet.Range("D10").Value = '2.70'
$q3Sheet.Range("E10").Value = '45.00'
$q3Sheet.Range("F10").Value = '2.15'
$q3Sheet.Range("G10").Value = '0.0580'
$q3Sheet.Range("H10").Value = 6
$q3Sheet.Range("A11").Value = 9
$q3Sheet.Range("B11").Value = '013326'
$q3Sheet.Range("C11").Value = '万家景气驱动混合A'
$q3Sheet.Range("D11").Value = '2.03'
$q3Sheet.Range("E11").Value = '92.93'
$q3Sheet.Range("F11").Value = '2.79'
$q3Sheet.Range("G11").Value = '0.0566'
$q3Sheet.Range("H11").Value = 10
$q3Sheet.Range("A12").Value = 10
$q3Sheet.Range("B12").Value = '519176'
$q3Sheet.Range("C12").Value = '浦银安盛消费升级混合C'
$q3Sheet.Range("D12").Value = '1.07'
$q3Sheet.Range("E12").Value = '90.82'
$q3Sheet.Range("F12").Value = '4.19'
$q3Sheet.Range("G12").Value = '0.0448'
$q3Sheet.Range("H12").Value = 10
$q3Sheet.Range("A13").Value = 11
$q3Sheet.Range("B13").Value = '015384'
$q3Sheet.Range("C13").Value = '万家瑞隆混合C'
$q3Sheet.Range("D13").Value = '1.52'
$q3Sheet.Range("E13").Value = '93.96'
$q3Sheet.Range("F13").Value = '2.85'
$q3Sheet.Range("G13").Value = '0.0433'
$q3Sheet.Range("H13").Value = 8
$q3Sheet.Range("A14").Value = 12
$q3Sheet.Range("B14").Value = '007735'
$q3Sheet.Range("C14").Value = '金鹰民安回报一年定期开放混合C'
$q3Sheet.Range("D14").Value = '2.59'
$q3Sheet.Range("E14").Value = '22.96'
$q3Sheet.Range("F14").Value = '1.55'
$q3Sheet.Range("G14").Value = '0.0401'
$q3Sheet.Range("H14").Value = 4
$q3Sheet.Range("A15").Value = 13
$q3Sheet.Range("B15").Value = '011351'
$q3Sheet.Range("C15").Value = '金鹰年年邮益一年持有期混合A'
$q3Sheet.Range("D15").Value = '3.43'
$q3Sheet.Range("E15").Value = '34.33'
$q3Sheet.Range("F15").Value = '1.12'
$q3Sheet.Range("G15").Value = '0.0384'
$q3Sheet.Range("H15").Value = 3
$q3Sheet.Range("A16").Value = 14
$q3Sheet.Range("B16").Value = '004819'
$q3Sheet.Range("C16").Value = '国寿安保目标策略灵活配置混合C'
$q3Sheet.Range("D16").Value = '1.73'
$q3Sheet.Range("E16").Value = '45.00'
$q3Sheet.Range("F16").Value = '2.15'
$q3Sheet.Range("G16").Value = '0.0372'
$q3Sheet.Range("H16").Value = 6
$q3Sheet.Range("A17").Value = 15
$q3Sheet.Range("B17").Value = '004265'
$q3Sheet.Range("C17").Value = '金鹰民丰回报定期开放混合'
$q3Sheet.Range("D17").Value = '4.53'
$q3Sheet.Range("E17").Value = '26.47'
$q3Sheet.Range("F17").Value = '0.82'
$q3Sheet.Range("G17").Value = '0.0371'
$q3Sheet.Range("H17").Value = 5
$q3Sheet.Range("A18").Value = 16
$q3Sheet.Range("B18").Value = '010853'
$q3Sheet.Range("C18").Value = '中欧内需成长混合C'
$q3Sheet.Range("D18").Value = '0.54'
$q3Sheet.Range("E18").Value = '90.11'
$q3Sheet.Range("F18").Value = '6.40'
$q3Sheet.Range("G18").Value = '0.0346'
$q3Sheet.Range("H18").Value = 5
$q3Sheet.Range("A19").Value = 17
$q3Sheet.Range("B19").Value = '008491'
$q3Sheet.Range("C19").Value = '万家周期优势企业混合A'
$q3Sheet.Range("D19").Value = '0.61'
$q3Sheet.Range("E19").Value = '93.50'
$q3Sheet.Range("F19").Value = '2.64'
$q3Sheet.Range("G19").Value = '0.0161'
$q3Sheet.Range("H19").Value = 9
$q3Sheet.Range("A20").Value = 18
$q3Sheet.Range("B20").Value = '013327'
$q3Sheet.Range("C20").Value = '万家景气驱动混合C'
$q3Sheet.Range("D20").Value = '0.32'
$q3Sheet.Range("E20").Value = '92.93'
$q3Sheet.Range("F20").Value = '2.79'
$q3Sheet.Range("G20").Value = '0.0089'
$q3Sheet.Range("H20").Value = 10
$q3Sheet.Range("A21").Value = 19
$q3Sheet.Range("B21").Value = '008492'
$q3Sheet.Range("C21").Value = '万家周期优势企业混合C'
$q3Sheet.Range("D21").Value = '0.14'
$q3Sheet.Range("E21").Value = '93.50'
$q3Sheet.Range("F21").Value = '2.64'
$q3Sheet.Range("G21").Value = '0.0037'
$q3Sheet.Range("H21").Value = 9
$q3Sheet.Range("A22").Value = 20
$q3Sheet.Range("B22").Value = '011352'
$q3Sheet.Range("C22").Value = '金鹰年年邮益一年持有期混合C'
$q3Sheet.Range("D22").Value = '0.27'
$q3Sheet.Range("E22").Value = '34.33'
$q3Sheet.Range("F22").Value = '1.12'
$q3Sheet.Range("G22").Value = '0.0030'
$q3Sheet.Range("H22").Value = 3

# The cloned sheet had 22 data rows (one more than Q3 needs); drop
# the now-unused trailing row so the used range matches (A1:H22)
$q3Sheet.Range("A23:H23").Clear()

# ---------------------------------------------------------------
# 3) Restore the active tab to the last sheet ("2021-Q1"), which
#    is where it was before our edits (Copy() made the new sheet
#    active as a side effect).
# ---------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()
